$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Team Members: reorder names, put Jim Cloud first; the _GoBack bookmark
#    (an artifact Word leaves at the point of the most recent edit) now sits
#    between "Jim Cloud, " and "Heather Mueller, Richard Xu".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Heather Mueller, Richard Xu, Jim Cloud", $true, $false, $false, $false,
    $false, $true, 1, $false, "Jim Cloud, Heather Mueller, Richard Xu", 2) | Out-Null

# Word drops the (hidden) _GoBack bookmark at the last edited spot, which is
# right before "Heather" now that the names have been reordered.
$goBackFind = $d.Content
$goBackFind.Find.Execute("Heather Mueller, Richard Xu", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackPoint = $d.Range($goBackFind.Start, $goBackFind.Start)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null

# ---------------------------------------------------------------------------
# 2) Description: drop the "ect based on social media and population
#    metrics" tail and tighten the wording.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Analyze shelter animal statistics for trends in adoption, euthanasia, return to owner, ect based on social media and population metrics.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Analyze shelter animal statistics for trends in adoption, euthanasia, and return to owner", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Data: no longer sourcing from the US Census API or social media -
#    remove those two data-source lines entirely.
# ---------------------------------------------------------------------------
$startRange = $d.Paragraphs.Item(9).Range.Start
$endRange = $d.Paragraphs.Item(10).Range.End
$d.Range($startRange, $endRange).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 4) Research questions: swap the social-media / population questions for
#    the new pet-type and location questions, tighten the Illinois-trends
#    question, and add a fourth question about listing a business email.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "How does social media outlets for shelters impact shelter impact outcomes of the      shelter animals?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Differences between adoption, euthanasia, and return to owner rates between different kinds of pets.", 2) | Out-Null

$d.Content.Find.Execute(
    "How does population impact the outcomes of the animals in shelters?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The effect of location on volume of adoption, euthanasia, and return to owner rates.  Which counties have the highest/lower amount, does this have correlation to population size or income levels?", 2) | Out-Null

$d.Content.Find.Execute(
    "What are the Illinois trends for intake, adoption, return to owner, and euthanasia over time.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Illinois trends for intake, adoption, return to owner, and euthanasia over time.  ", 2) | Out-Null

$thirdQuestion = $d.Paragraphs.Item(14)
$newQuestion = $thirdQuestion.Range.InsertParagraphAfter()
$fourthQuestion = $d.Paragraphs.Item(15)
$fourthQuestion.Range.Text = "Does listing a business email have any effects of rates?"
